# CCC19 Derived Variables Spreadsheet
# Add a new derived variable "steroids_c19" (Rx04a): ever/never treatment of
# COVID-19 with steroids, not taking dosing into account. Because the existing
# "Rx04" variable ("steroids") only ever captured HIGH-dose steroid exposure,
# rename it to "steroids_hd_c19" to disambiguate it from the new variable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Rx04 row ("High-dose steroids as COVID-19 treatment ever") currently sits on
# row 87, with variable name "steroids" -- rename to "steroids_hd_c19".
$ws.Range("B87").Value = "steroids_hd_c19"

# Insert a brand-new row right below it (row 88) for the new variable, and
# grow the table/list-object so the new row is part of Table1.
$ws.Rows.Item(88).Insert()
$lo.Resize($ws.Range("A1:E118"))

$ws.Range("A88").Value = "Rx04a"
$ws.Range("B88").Value = "steroids_c19"
$ws.Range("C88").Value = "Treatments"
$ws.Range("D88").Value = "Steroids as COVID-19 treatment ever"

# Match the author's on-screen selection after the edit.
$ws.Range("D88").Select()
